$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "43.831.67"
Set-TextValue 2 5 "  -0.79%  "
Set-TextValue 3 4 "2.232.69"
Set-TextValue 3 5 "  -2.84%  "
Set-TextValue 4 4 "1.01"
Set-TextValue 4 5 "  +0.69%  "
Set-TextValue 5 4 "0.650"
Set-TextValue 5 5 "  +3.49%  "
Set-TextValue 6 4 "229.65"
Set-TextValue 6 5 "  -0.64%  "
Set-TextValue 7 4 "62.32"
Set-TextValue 7 5 "  +1.84%  "
Set-TextValue 8 5 "  +0.03%  "
Set-TextValue 9 4 "0.445"
Set-TextValue 9 5 "  +4.14%  "
Set-TextValue 10 4 "0.0956"
Set-TextValue 10 5 "  +0.91%  "
Set-TextValue 11 4 "56.80"
Set-TextValue 11 5 "  -2.18%  "
Set-TextValue 12 4 "26.37"
Set-TextValue 12 5 "  +7.58%  "
Set-TextValue 13 5 "  +1.26%  "
Set-TextValue 14 4 "2.567.73"
Set-TextValue 14 5 "  -2.64%  "
Set-TextValue 15 4 "15.33"
Set-TextValue 15 5 "  -2.71%  "
Set-TextValue 16 4 "6.06"
Set-TextValue 16 5 "  +2.25%  "
Set-TextValue 17 5 "  +0.63%  "
Set-TextValue 18 4 "2.229.55"
Set-TextValue 18 5 "  -2.98%  "
Set-TextValue 19 4 "43.649.18"
Set-TextValue 19 5 "  -0.92%  "
Set-TextValue 20 4 "0.0₃0979"
Set-TextValue 20 5 "  +3.12%  "
Set-TextValue 21 4 "72.40"
Set-TextValue 21 5 "  -2.42%  "
Set-TextValue 22 4 "5.99"
Set-TextValue 22 5 "  -4.33%  "
Set-TextValue 23 4 "247.60"
Set-TextValue 23 5 "  -2.68%  "
Set-TextValue 25 5 "  -5.83%  "
Set-TextValue 26 2 "WEMIXToken"
Set-TextValue 26 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 26 4 "3.38"
Set-TextValue 26 5 "  +23.18%  "
Set-TextValue 27 2 "Toncoin"
Set-TextValue 27 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue 27 4 "2.23"
Set-TextValue 27 5 "  -5.47%  "
Set-TextValue 28 4 "9.77"
Set-TextValue 28 5 "  -1.43%  "
Set-TextValue 29 4 "170.22"
Set-TextValue 29 5 "  -0.35%  "
Set-TextValue 30 5 "  -1.90%  "
Set-TextValue 31 4 "20.69"
Set-TextValue 31 5 "  +0.23%  "
Set-TextValue 32 4 "1.39"
Set-TextValue 32 5 "  -2.77%  "
Set-TextValue 33 5 "  +3.10%  "
Set-TextValue 34 5 "  +4.87%  "
Set-TextValue 35 4 "4.72"
Set-TextValue 35 5 "  -1.14%  "
Set-TextValue 36 4 "4.87"
Set-TextValue 36 5 "  -4.20%  "
Set-TextValue 37 5 "  -1.15%  "
Set-TextValue 38 4 "6.36"
Set-TextValue 38 5 "  -2.84%  "
Set-TextValue 39 5 "  -5.93%  "
Set-TextValue 40 4 "0.0256"
Set-TextValue 40 5 "  +2.19%  "
Set-TextValue 41 5 "  -0.11%  "
Set-TextValue 42 5 "  -3.33%  "
Set-TextValue 43 4 "8.20"
Set-TextValue 43 5 "  -5.66%  "
Set-TextValue 44 4 "16.93"
Set-TextValue 44 5 "  -1.48%  "
Set-TextValue 45 4 "0.0948"
Set-TextValue 45 5 "  -2.56%  "
Set-TextValue 46 4 "96.44"
Set-TextValue 46 5 "  -2.52%  "
Set-TextValue 47 5 "  -2.88%  "
Set-TextValue 48 4 "4.34"
Set-TextValue 48 5 "  -1.28%  "
Set-TextValue 49 4 "2.30"
Set-TextValue 49 5 "  +1.61%  "
Set-TextValue 50 4 "1.424.44"
Set-TextValue 50 5 "  -4.12%  "
Set-TextValue 51 5 "  +1.76%  "
